$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26, 27, 28 get their data (columns A, B, E, F, G, H, Q, R) cyclically
# rotated: old row26 -> new row27, old row27 -> new row28, old row28 -> new row26.
$ws.Range("A26").Value = 111541116
$ws.Range("B26").Value = 108219
$ws.Range("E26").Value = 219711
$ws.Range("F26").Value = "Sårläka"
$ws.Range("G26").Value = "Sanicula europaea"
$ws.Range("H26").Value = "L."
$ws.Range("Q26").Value = 693830.7552326696
$ws.Range("R26").Value = 6552178.401404973

$ws.Range("A27").Value = 111541130
$ws.Range("B27").Value = 98535
$ws.Range("E27").Value = 222498
$ws.Range("F27").Value = "Blåsippa"
$ws.Range("G27").Value = "Hepatica nobilis"
$ws.Range("H27").Value = "Schreb."
$ws.Range("Q27").Value = 693830.8333423812
$ws.Range("R27").Value = 6552176.860022029

$ws.Range("A28").Value = 111541117
$ws.Range("B28").Value = 108219
$ws.Range("E28").Value = 219711
$ws.Range("F28").Value = "Sårläka"
$ws.Range("G28").Value = "Sanicula europaea"
$ws.Range("H28").Value = "L."
$ws.Range("Q28").Value = 693809.5100469354
$ws.Range("R28").Value = 6552200.504896822
